$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the template error: A1 should hold the header label "Parameter_Name"
# (it previously contained the stray numeric value 3).
$ws.Range("A1").Value = "Parameter_Name"

# Fix link/selection on the page: make A2 the active/selected cell.
$ws.Range("A2").Select()
